# Auto update: 2025-11-29 03:23:53
# Update hedging/insurance analysis figures for PRU, UNH, MET, AIG rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Prudential Financial, Inc. (PRU)
$ws.Range("D2").Value = 108.44
$ws.Range("F2").Value = 4.15
$ws.Range("K2").Value = 59.2
$ws.Range("N2").Value = 85.96878041621773

# Row 3 - UnitedHealth Group Incorporated (UNH)
$ws.Range("D3").Value = 330.04
$ws.Range("F3").Value = 5.94
$ws.Range("K3").Value = 56.2
$ws.Range("N3").Value = 85.96878041621773

# Row 4 - MetLife, Inc. (MET)
$ws.Range("D4").Value = 76.61
$ws.Range("F4").Value = 3.64
$ws.Range("I4").Value = 46
$ws.Range("K4").Value = 53.2
$ws.Range("N4").Value = 85.96878041621773

# Row 5 - American International Group, Inc. (AIG)
$ws.Range("D5").Value = 76.34999999999999
$ws.Range("F5").Value = 0.87
$ws.Range("K5").Value = 51.8
$ws.Range("N5").Value = 85.96878041621773
